# Refresh the cryptos price/volume snapshot (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D holds prices formatted as plain text (e.g. "259.30"). Excel
# auto-converts such strings to numbers on assignment (dropping trailing
# zeros / turning "70.00" into 70), so force those specific cells to Text
# format first to preserve the exact published string.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '37.264.62'
$ws.Range("E2").Value = '  -0.22%  '
$ws.Range("D3").Value = '2.006.04'
$ws.Range("E3").Value = '  -1.41%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '259.30'
$ws.Range("E5").Value = '  +4.35%  '
$ws.Range("E6").Value = '  -1.96%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '56.49'
$ws.Range("E8").Value = '  -7.16%  '
$ws.Range("E9").Value = '  -3.75%  '
$ws.Range("D10").Value = '0.0770'
$ws.Range("E10").Value = '  -5.22%  '
$ws.Range("E11").Value = '  -3.13%  '
$ws.Range("D12").Value = '2.305.87'
$ws.Range("E12").Value = '  -1.26%  '
$ws.Range("D13").Value = '14.22'
$ws.Range("E13").Value = '  -7.35%  '
$ws.Range("D14").Value = '21.69'
$ws.Range("E14").Value = '  -3.35%  '
$ws.Range("E15").Value = '  -8.14%  '
$ws.Range("E16").Value = '  -6.31%  '
$ws.Range("D17").Value = '2.039.96'
$ws.Range("E17").Value = '  +0.25%  '
$ws.Range("D18").Value = '37.278.65'
$ws.Range("E18").Value = '  -0.09%  '
$ws.Range("D19").Value = '70.00'
$ws.Range("E19").Value = '  -1.19%  '
$ws.Range("D20").Value = '0.0₃0834'
$ws.Range("E20").Value = '  -4.03%  '
$ws.Range("D21").Value = '232.08'
$ws.Range("E21").Value = '  +0.46%  '
$ws.Range("E22").Value = '  -3.27%  '
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("D24").Value = '2.58'
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("E25").Value = '  -0.40%  '
$ws.Range("D26").Value = '164.57'
$ws.Range("E26").Value = '  +0.38%  '
$ws.Range("D27").Value = '8.93'
$ws.Range("E27").Value = '  -6.11%  '
$ws.Range("D28").Value = '19.53'
$ws.Range("E28").Value = '  -1.73%  '
$ws.Range("E29").Value = '  -6.35%  '
$ws.Range("D30").Value = '1.32'
$ws.Range("E30").Value = '  -4.95%  '
$ws.Range("E31").Value = '  -2.31%  '
$ws.Range("E32").Value = '  -5.57%  '
$ws.Range("D33").Value = '0.0639'
$ws.Range("E33").Value = '  -5.15%  '
$ws.Range("E34").Value = '  -1.75%  '
$ws.Range("D35").Value = '2.35'
$ws.Range("E35").Value = '  -6.56%  '
$ws.Range("E36").Value = '  +0.24%  '
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("E38").Value = '  -8.76%  '
$ws.Range("E39").Value = '  -0.91%  '
$ws.Range("E40").Value = '  +1.62%  '
$ws.Range("E41").Value = '  -1.67%  '
$ws.Range("E42").Value = '  -1.98%  '
$ws.Range("D43").Value = '0.0924'
$ws.Range("E43").Value = '  -6.31%  '
$ws.Range("D44").Value = '1.430.13'
$ws.Range("E44").Value = '  +3.15%  '
$ws.Range("D45").Value = '89.33'
$ws.Range("E45").Value = '  -3.95%  '
$ws.Range("D46").Value = '15.64'
$ws.Range("E46").Value = '  -8.97%  '
$ws.Range("E47").Value = '  -4.14%  '
$ws.Range("D48").Value = '2.92'
$ws.Range("E48").Value = '  +2.11%  '
$ws.Range("D49").Value = '6.99'
$ws.Range("E49").Value = '  -6.80%  '
$ws.Range("D50").Value = '2.198.40'
$ws.Range("E50").Value = '  -1.20%  '
$ws.Range("D51").Value = '1.94'
$ws.Range("E51").Value = '  -11.26%  '
